# TC_MISC_GA_WH_Yes_01 case flow completed
#
# 1) Update a handful of result cells on the "1099MISCdata" sheet (new test-run
#    timestamps / EIN / reference numbers) and append a small "reconciliation"
#    summary block under the existing data.
# 2) Add a brand-new "StateID-Validations" worksheet after "1099MISCdata" with
#    a small validation matrix for state-id formats.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "1099MISCdata" sheet updates
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("1099MISCdata")

$ws.Range("B3").Value  = "Test0305202183817"
$ws.Range("B4").Value  = "Test0305202183207"
$ws.Range("B5").Value  = "201204660"
$ws.Range("B9").Value  = "80-1825910"
$ws.Range("B16").Value = 236894931

$ws.Range("A22").Value = "Reconsillation details"
$ws.Range("A23").Value = "State name"
$ws.Range("A24").Value = "Form count"
$ws.Range("A25").Value = "Status"

$ws.Range("B23").Select()

# ---------------------------------------------------------------------------
# 2) New "StateID-Validations" sheet
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newWs = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newWs.Name = "StateID-Validations"

$newWs.Range("E1").Value = "State"
$newWs.Range("F1").Value = "GA"

$newWs.Range("A2").Value = 1
$newWs.Range("B2").Value = 123
$newWs.Range("C2").Value = "Less than 9 digits"

$newWs.Range("A3").Value = 2
$newWs.Range("B3").Value = 1234567891
$newWs.Range("C3").Value = "Greater than 9 digits"

$newWs.Range("A4").Value = 3
$newWs.Range("B4").Value = "123456789kk"
$newWs.Range("C4").Value = "Missed hyphen"

$newWs.Range("A5").Value = 4
$newWs.Range("B5").Value = "1234567-k2"
$newWs.Range("C5").Value = "Should have 7 digits and 2 Alphabets"

$newWs.Range("A6").Value = 5
$newWs.Range("B6").Value = "SS34567-KJ"
$newWs.Range("C6").Value = "Should have 7 digits,hyphen and 2 Alphabets"

$newWs.Range("A7").Value = 6
$newWs.Range("B7").Value = "1234567-AA"
$newWs.Range("C7").Value = "Correct format value has 7 digits and 2 Alphabets"

$newWs.Range("A8").Value = 7
$newWs.Range("B8").Value = 123456789
$newWs.Range("C8").Value = "Correct format value has 9 digits"

$newWs.Range("B8").Select()
